{"js": "// The requirement paragraph \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito\n// fraco)\" is followed by four paragraphs that are no longer wanted:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. an empty paragraph\n//   4. an empty paragraph with pageBreakBefore\n// Locate the anchor paragraph by its text, walk forward four paragraphs,\n// then delete them (removing the later ones first so the earlier\n// `Paragraph` objects stay valid).\n\nconst body = context.document.body;\nconst searchResults = body.search(\n  \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Anchor paragraph \"LOQ4086: ...\" not found.');\n}\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n\nconst toRemove1 = anchorParagraph.getNext();\nconst toRemove2 = toRemove1.getNext();\nconst toRemove3 = toRemove2.getNext();\nconst toRemove4 = toRemove3.getNext();\nawait context.sync();\n\n// Delete in reverse order so earlier paragraph references remain valid.\ntoRemove4.delete();\ntoRemove3.delete();\ntoRemove2.delete();\ntoRemove1.delete();\n\nawait context.sync();\n", "ps1": "# The requirement paragraph \"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito\n# fraco)\" is followed by four paragraphs that should be removed:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. an empty paragraph\n#   4. an empty paragraph with PageBreakBefore\n# Find the anchor paragraph by its text, walk forward four paragraphs via\n# Paragraph.Next(), then delete their ranges starting from the last one so\n# earlier paragraph references stay valid.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"LOQ4086: Opera\u00e7\u00f5es Unit\u00e1rias II (Requisito fraco)\")\n\nif (-not $found) {\n    throw 'Anchor paragraph \"LOQ4086: ...\" not found.'\n}\n\n$anchorParagraph = $searchRange.Paragraphs(1)\n\n$toRemove1 = $anchorParagraph.Next()\n$toRemove2 = $toRemove1.Next()\n$toRemove3 = $toRemove2.Next()\n$toRemove4 = $toRemove3.Next()\n\n# Delete in reverse order so earlier paragraph references remain valid.\n$toRemove4.Range.Delete()\n$toRemove3.Range.Delete()\n$toRemove2.Range.Delete()\n$toRemove1.Range.Delete()\n"}
